# trend_2019_weekly.xlsx - weekly trend sheet update ("now it can use")
# - bold the header row
# - rename the existing KPI row to the "overall" qualifying-rate wording
#   and drop its old highlight fill
# - add three new KPI rows (shape / surface / composition-performance),
#   each owned by 王宇阳
# - move the active selection to E10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row (row 1): bold the whole row ----
$ws.Rows.Item(1).Font.Bold = $true

# ---- row 2: rename the KPI label, drop the old highlight fill ----
$ws.Range("A2").Value = "1580高牌号硅钢总体一次投料合格率"
$ws.Range("A2:H2").ClearFormats()
$ws.Range("B2:D2").VerticalAlignment = -4107

# ---- new KPI rows, all owned by 王宇阳 ----
$ws.Range("A3").Value = "1580高牌号硅钢板形一次投料合格率"
$ws.Range("B3").Value = "王宇阳"
$ws.Range("C3").Value = 92.76
$ws.Range("D3").Value = 94.77
$ws.Range("A3:D3").VerticalAlignment = -4107

$ws.Range("A4").Value = "1580高牌号硅钢表面一次投料合格率"
$ws.Range("B4").Value = "王宇阳"
$ws.Range("C4").Value = 88.21
$ws.Range("D4").Value = 91.65
$ws.Range("A4:D4").VerticalAlignment = -4107

$ws.Range("A5").Value = "1580高牌号硅钢成份性能一次投料合格率"
$ws.Range("B5").Value = "王宇阳"
$ws.Range("C5").Value = 71.12
$ws.Range("D5").Value = 77.85
$ws.Range("A5:D5").VerticalAlignment = -4107

# ---- selection ----
$ws.Range("E10").Select()
